# Update column F (dSF) values on Sheet1 for the rows whose "dSF" figure
# was repulled/recalculated, per the commit "repull data, push all data,
# mean calculation".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @{
    3  = 2
    6  = -1
    15 = 2
    18 = 2
    25 = 1
    41 = 0
    42 = 1
    47 = -3
    50 = 0
    52 = 0
    57 = -2
    58 = -2
    65 = -1
    68 = 4
    69 = 1
    71 = -1
    75 = -1
    77 = 0
    82 = -4
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
